$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    @(2, -0.3014486003302357, -0.2836103601562896),
    @(3, 0.2208790396742221, 0.2310613204452885),
    @(4, 0.1629816387700927, 0.1755306681328835),
    @(5, -0.05012634322047164, -0.04056665752737425),
    @(6, 0.1561968451442523, 0.17212812949264),
    @(7, -0.4370487854644528, -0.4265721493985016),
    @(8, -0.2825953894762833, -0.2609341954564389),
    @(9, -0.405010385582875, -0.3946392881059867),
    @(10, 0.3284768669770559, 0.345795008925863),
    @(11, -0.2029198406968479, -0.1949870406157501),
    @(12, -0.0829859070944069, -0.05739727746732613),
    @(13, 0.001860657265781648, 0.002311613583817898),
    @(14, 0.03690761557091227, 0.07681899866205064),
    @(15, -0.08004324378204768, -0.03093631638004116),
    @(16, 0.265441533111786, 0.3314487172340432),
    @(17, 0.5784939741767463, 0.6110262923049181),
    @(18, 0.0572461162532695, 0.0388935416415849),
    @(19, 0.3993796273122018, 0.4084822820252282),
    @(20, 0.2392378157152579, 0.2983112640309355),
    @(21, 0.4409231861868707, 0.5079155085024114),
    @(22, 0.3489993189955894, 0.3856127404050427),
    @(23, -0.06826607600014872, -0.03634721152122115),
    @(24, 4.542957436487145, 4.570757146750512),
    @(25, 0.5297536296279391, 0.4912331328330304),
    @(26, 0.389171018658315, 0.3684126309347704),
    @(27, 0.3409303774337827, 0.3124607669649356),
    @(28, 1.103567232980138, 1.071769410976996),
    @(29, 5.852343364666749, 5.402542441936599),
    @(30, 1.034456036008815, 0.9825880294858924),
    @(31, -0.1305243195791115, -0.1954689676156292),
    @(32, 0.8241433089858108, 0.7863655111848381),
    @(33, 0.9213854168362843, 0.8944894956422242),
    @(34, -0.6037134917095873, -0.63624695351572),
    @(35, 0.8322579948917794, 0.8223535213348641),
    @(36, 0.7860084475941813, 0.7698895913431909),
    @(37, 0.7681381328274761, 0.7474376185359266),
    @(38, 0.7602410526304069, 0.7374706735306866),
    @(39, 0.5828317397742391, 0.5804290068628625),
    @(40, 0.7541647801922936, 0.7528309752724952),
    @(41, 0.5738368562037698, 0.5660837628946727),
    @(42, 0.7225480528620418, 0.6931587004027192),
    @(43, 0.7340831943415702, 0.7184183177175822),
    @(44, 0.674460994729557, 0.6672824845645405),
    @(45, 0.6605212839239608, 0.6433336040147427),
    @(46, -1.252385680772581, -1.256773037854898),
    @(47, -0.9696040711523473, -0.9747718144687532),
    @(48, -0.8624287680479477, -0.8687396927531337),
    @(49, -0.6314397035640043, -0.6349963409332019),
    @(50, -0.04881614723499951, -0.04925852735632676),
    @(51, -0.8513797302704933, -0.8554818286696787),
    @(52, -0.8513797302704933, -0.8554818286696787),
    @(53, -1.078713386057442, -1.093424245395886),
    @(54, -0.1865630227146239, -0.1851023459123787),
    @(55, -0.9898708584296871, -0.9936145961238633),
    @(56, -0.8922206499949255, -0.8864877234731758),
    @(57, -0.9656975326397896, -0.9487388770193534),
    @(58, -1.167249448836082, -1.136821590379815),
    @(59, -0.870412731334352, -0.8529240316513477),
    @(60, -0.5226154975282774, -0.4968076412149591),
    @(61, 0.3666051981124333, 0.3696770819872654),
    @(62, -1.256446171571648, -1.233318256826487),
    @(63, -0.7646879888823394, -0.72945411379562),
    @(64, -0.9016649760301245, -0.893155543773116),
    @(65, -0.1412741075018785, -0.1157307761025997),
    @(66, -0.8349009314063618, -0.8041248071301309),
    @(67, -0.8376161498186951, -0.7953707832833182)
)

foreach ($row in $newValues) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value2 = $row[1]
    $ws.Cells.Item($r, 3).Value2 = $row[2]
}
